# Update the "取得日時" (acquired timestamp) column A for rows 2-8
# on the "ランサーズ" sheet from 2025-10-26 12:33:22 to 2025-10-26 12:44:08.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-26 12:44:08"

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
